# Swap the data between row 10 and row 11 for columns A, B, E, F, G, H.
# This corresponds to the diff, where the bird record for "Rödvingetrast"
# (Turdus iliacus) and the record for "Björktrast" (Turdus pilaris) swap
# row positions (row 10 <-> row 11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$cols = @("A", "B", "E", "F", "G", "H")

foreach ($col in $cols) {
    $cell10 = $ws.Range($col + "10")
    $cell11 = $ws.Range($col + "11")

    $val10 = $cell10.Value()
    $val11 = $cell11.Value()

    $cell10.Value = $val11
    $cell11.Value = $val10
}
